$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# Sheet "IUCN 2019 Redlist": replace the summary table with the expanded
# whale population dataset (adds CCE / Christensen(2006) / historical
# columns plus four right-whale / bowhead species rows).
# ----------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("IUCN 2019 Redlist")
$ws1.Activate()

# Wipe the old 7-column table (values + formatting) before laying out
# the new 14-column one.
$ws1.Cells.Clear()

$data = New-Object 'object[,]' 12,14
$data[0,0] = "Species"
$data[0,1] = "Population estimate (IUCN 2019)"
$data[0,2] = "CCE population estimate"
$data[0,3] = "Southern hemisphere population estimate (Christensen 2006)"
$data[0,4] = "Population estimate (Christensen 2006)"
$data[0,5] = "Population low estimate (Christensen 2006)"
$data[0,6] = "Population high estimate (Christensen 2006)"
$data[0,7] = "Historical estimate"
$data[0,8] = "Historical low estimate"
$data[0,9] = "Historical high estimate"
$data[0,10] = "Southern hemisphere historic estimate (Christensen 2006)"
$data[0,11] = "Number removed by 20th century whaling (N. Hemisphere)"
$data[0,12] = "Number removed by 20th century whaling (S. Hemisphere)"
$data[0,13] = "Total removed"
$data[1,0] = "Balaenoptera musculus"
$data[1,1] = 10000
$data[1,2] = 1647
$data[1,3] = 1180
$data[1,4] = 4727
$data[1,5] = 3378
$data[1,6] = 6181
$data[1,7] = 340280
$data[1,8] = 308510
$data[1,9] = 376120
$data[1,10] = 327000
$data[1,11] = 15537
$data[1,12] = 363648
$data[1,13] = $null
$data[2,0] = "Balaenoptera physalus"
$data[2,1] = 100000
$data[2,2] = 9029
$data[2,3] = 55700
$data[2,4] = 109600
$data[2,5] = 72200
$data[2,6] = 161200
$data[2,7] = 762400
$data[2,8] = 573500
$data[2,9] = 936000
$data[2,10] = 625000
$data[2,11] = 147607
$data[2,12] = 726461
$data[2,13] = $null
$data[3,0] = "Megaptera novaeangliae"
$data[3,1] = 84000
$data[3,2] = 1918
$data[3,3] = 22500
$data[3,4] = 42070
$data[3,5] = 31510
$data[3,6] = 59000
$data[3,7] = 231700
$data[3,8] = 154500
$data[3,9] = 285400
$data[3,10] = 199000
$data[3,11] = 33585
$data[3,12] = 215848
$data[3,13] = $null
$data[4,0] = "Balaenoptera borealis"
$data[4,1] = 50000
$data[4,2] = 519
$data[4,3] = 6990
$data[4,4] = 49090
$data[4,5] = 27780
$data[4,6] = 75740
$data[4,7] = 246000
$data[4,8] = 219020
$data[4,9] = 294400
$data[4,10] = 167000
$data[4,11] = 86951
$data[4,12] = 204589
$data[4,13] = $null
$data[5,0] = "Balaenoptera edeni"
$data[5,1] = 80000
$data[5,2] = $null
$data[5,3] = 91300
$data[5,4] = 132400
$data[5,5] = 97600
$data[5,6] = 176500
$data[5,7] = 146300
$data[5,8] = 111600
$data[5,9] = 190800
$data[5,10] = 94100
$data[5,11] = 14049
$data[5,12] = 7913
$data[5,13] = $null
$data[6,0] = "Balaenoptera acutorostrata"
$data[6,1] = 200000
$data[6,2] = 636
$data[6,3] = $null
$data[6,4] = 188900
$data[6,5] = 141900
$data[6,6] = 251400
$data[6,7] = 258000
$data[6,8] = 195700
$data[6,9] = 344300
$data[6,10] = $null
$data[6,11] = 166692
$data[6,12] = $null
$data[6,13] = $null
$data[7,0] = "Balaenoptera bonaerensis"
$data[7,1] = 515000
$data[7,2] = $null
$data[7,3] = 318000
$data[7,4] = 318000
$data[7,5] = 250000
$data[7,6] = 404000
$data[7,7] = 379000
$data[7,8] = 300000
$data[7,9] = 478000
$data[7,10] = 379000
$data[7,11] = $null
$data[7,12] = 117213
$data[7,13] = $null
$data[8,0] = "Eubalaena glacialis"
$data[8,1] = $null
$data[8,2] = $null
$data[8,3] = $null
$data[8,4] = 6740
$data[8,5] = 4580
$data[8,6] = 11100
$data[8,7] = 14100
$data[8,8] = 10100
$data[8,9] = 27800
$data[8,10] = $null
$data[8,11] = $null
$data[8,12] = 141
$data[8,13] = $null
$data[9,0] = "Eubalaena japonica"
$data[9,1] = $null
$data[9,2] = $null
$data[9,3] = $null
$data[9,4] = 368
$data[9,5] = 257
$data[9,6] = 469
$data[9,7] = 14100
$data[9,8] = 10100
$data[9,9] = 27800
$data[9,10] = $null
$data[9,11] = 967
$data[9,12] = $null
$data[9,13] = $null
$data[10,0] = "Eubalaena australis"
$data[10,1] = $null
$data[10,2] = $null
$data[10,3] = 6740
$data[10,4] = 368
$data[10,5] = 257
$data[10,6] = 469
$data[10,7] = 86100
$data[10,8] = 73400
$data[10,9] = 98300
$data[10,10] = 86100
$data[10,11] = $null
$data[10,12] = 4452
$data[10,13] = $null
$data[11,0] = "Balaena mysticetus"
$data[11,1] = 10000
$data[11,2] = $null
$data[11,3] = $null
$data[11,4] = 9450
$data[11,5] = 7500
$data[11,6] = 10800
$data[11,7] = 89000
$data[11,8] = 67000
$data[11,9] = 114000
$data[11,10] = $null
$data[11,11] = $null
$data[11,12] = $null
$data[11,13] = $null


$ws1.Range("A1:N12").Value = $data

# N column: running total (N2:N9 as a live formula, N10:N11 as the
# literal values the source workbook stores, N12 left blank).
$ws1.Range("N2").Formula = "=L2+M2"
$ws1.Range("N3:N9").Formula = "=L3+M3"
$ws1.Range("N10").Value = 967
$ws1.Range("N11").Value = 4452

# Header row formatting (bold, already the sheet's row-1 style) plus the
# one-off bold-black-Calibri styling on the "Historical high estimate"
# header cell (J1).
$ws1.Range("A1:N1").Font.Bold = $true
$ws1.Range("J1").Font.Name = "Calibri"
$ws1.Range("J1").Font.Bold = $true
$ws1.Range("J1").Font.Color = 0

# Column widths (approximate match to the authored layout).
$ws1.Columns.Item(1).ColumnWidth = 23.1666666
$ws1.Columns.Item(2).ColumnWidth = 27.8307291
$ws1.Columns.Item(3).ColumnWidth = 20.6666666
$ws1.Columns.Item(4).ColumnWidth = 51.9986979
$ws1.Columns.Item(5).ColumnWidth = 35.4986979
$ws1.Columns.Item(6).ColumnWidth = 35.4986979
$ws1.Columns.Item(7).ColumnWidth = 35.4986979
$ws1.Columns.Item(8).ColumnWidth = 16.1666666
$ws1.Columns.Item(9).ColumnWidth = 19.6666666
$ws1.Columns.Item(10).ColumnWidth = 20.1666666
$ws1.Columns.Item(11).ColumnWidth = 49.1666666
$ws1.Columns.Item(12).ColumnWidth = 50.4986979
$ws1.Columns.Item(13).ColumnWidth = 50.1666666

# Freeze the species column and leave the selection / scroll position
# where the author left it.
$ws1.Range("B1").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws1.Range("K13").Select()

# ----------------------------------------------------------------------
# Sheet "Resources": bold the "Christensen 2006" citation and append the
# new "Rocha et al. 2014" / "Whaling numbers" reference row.
# ----------------------------------------------------------------------
$ws6 = $wb.Worksheets.Item("Resources")
$ws6.Range("A4").Font.Bold = $true
$ws6.Range("A6").Value = "Rocha et al. 2014"
$ws6.Range("C6").Value = "Whaling numbers"
